$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.987.09"
$ws.Range("E2").Value = "  -0.37%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.741.19"
$ws.Range("E3").Value = "  +0.01%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.93"
$ws.Range("E5").Value = "  +1.50%  "

$ws.Range("E6").Value = "  -0.04%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5042"
$ws.Range("E7").Value = "  -4.65%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2740"
$ws.Range("E8").Value = "  -1.81%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06175"
$ws.Range("E9").Value = "  -0.01%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.748.29"
$ws.Range("E10").Value = "  +0.18%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07246"
$ws.Range("E11").Value = "  +0.79%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.6508"
$ws.Range("E12").Value = "  -0.47%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.10"
$ws.Range("E13").Value = "  -1.16%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.672"
$ws.Range("E14").Value = "  +0.84%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.41"
$ws.Range("E15").Value = "  -0.50%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("E16").Value = "  +0.04%  "

$ws.Range("E17").Value = "  -0.07%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.000.65"
$ws.Range("E18").Value = "  -0.36%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.88"
$ws.Range("E19").Value = "  +0.16%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000006879"
$ws.Range("E20").Value = "  +1.55%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.973.29"
$ws.Range("E21").Value = "  +0.18%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.459"
$ws.Range("E22").Value = "  +0.44%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.705"
$ws.Range("E23").Value = "  -0.15%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.360"
$ws.Range("E24").Value = "  +2.07%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "135.69"
$ws.Range("E25").Value = "  -3.39%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.505"
$ws.Range("E26").Value = "  -0.76%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.25"
$ws.Range("E27").Value = "  -0.28%  "

$ws.Range("E28").Value = "  -0.19%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "105.35"
$ws.Range("E29").Value = "  -0.36%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.884"
$ws.Range("E30").Value = "  +0.69%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08181"
$ws.Range("E31").Value = "  -3.24%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.641"
$ws.Range("E32").Value = "  -0.62%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04660"

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.660"
$ws.Range("E34").Value = "  -0.26%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9941"
$ws.Range("E35").Value = "  -0.13%  "

$ws.Range("B36").Value = "MXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.801"
$ws.Range("E36").Value = "  +3.66%  "

$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6078"
$ws.Range("E37").Value = "  -3.13%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01621"
$ws.Range("E38").Value = "  +0.38%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.922"
$ws.Range("E39").Value = "  -0.48%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.001"
$ws.Range("E40").Value = "  -0.01%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "100.55"
$ws.Range("E41").Value = "  +0.91%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.3906"
$ws.Range("E42").Value = "  +0.09%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7635"
$ws.Range("E43").Value = "  +1.25%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.000"
$ws.Range("E44").Value = "  +0.16%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1161"
$ws.Range("E45").Value = "  +1.06%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.296"
$ws.Range("E46").Value = "  +0.08%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.60"
$ws.Range("E47").Value = "  +1.26%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05314"
$ws.Range("E48").Value = "  -0.28%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "30.60"
$ws.Range("E49").Value = "  -0.72%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.624"
$ws.Range("E50").Value = "  +0.38%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3451"
$ws.Range("E51").Value = "  -0.71%  "
